# Apply the VO2max trend CodeSystem metadata fixes described by the diff:
#  - B7  (Experimental) : empty -> "false"
#  - B8  (Date)          : "2025-11-28T14:35:57+00:00" -> "2025-11-30T13:08:37+00:00"
#  - B17 (Description)   : empty -> "Codes for VO2max trend direction over time"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- B7: Experimental = "false" ---------------------------------------
# A literal Value2 assignment of "false"/"true" is auto-coerced to a
# boolean by the engine (as real Excel does for unformatted text typed
# directly into a cell), which is not what we want here: the target
# workbook stores this as literal text "false". Build it as a formula
# that evaluates to the text string, then convert the formula to a
# plain value via copy / paste-special-values so the stored cell ends
# up as ordinary text (not a logical value, and without mutating the
# cell's style the way an apostrophe-quoted literal would).
$quote = [char]34
$falseFormula = "=CONCATENATE(" + $quote + "fal" + $quote + "," + $quote + "se" + $quote + ")"
$ws.Range("B7").Formula = $falseFormula
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)   # xlPasteValues

# --- B8: Date -----------------------------------------------------------
$ws.Range("B8").Value2 = "2025-11-30T13:08:37+00:00"

# --- B17: Description -----------------------------------------------
$ws.Range("B17").Value2 = "Codes for VO2max trend direction over time"
